$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-18 01:14:14"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
